# Auto-generated edit script to apply the diff changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.419.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.853.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.21"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6925"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3058"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07650"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.48"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07743"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.138"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.846.12"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6927"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.83"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.314"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.434.82"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.098.75"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.34"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.71"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.654"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9999"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1474"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.954"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.94"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.533"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.248"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.137"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.204"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05210"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7733"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.872"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.145"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.684"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.322.94"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +8.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01866"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.718"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9425"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.01"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.786"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9995"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.717"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.998.52"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5225"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.31%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.779"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000123"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "62.91"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05954"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.00%  "
